# Restored from revision of admin on 09/30/2020 07:17:57 AM.TEST Author: admin. Type: SAVE.
# The rule table's row 10 ("R30") "Integer min" threshold changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
